$d = $word.ActiveDocument

function Expand-SimpleField($doc, $instr) {
    # Locate the fldSimple field whose instruction text matches $instr and
    # replace it with the expanded begin/instrText/separate/end run sequence,
    # exactly as Word does when a simple field is turned into a complex one.
    $target = $null
    for ($i = 1; $i -le $doc.Fields.Count; $i++) {
        $f = $doc.Fields.Item($i)
        if ($f.Code.Text -eq $instr) {
            $target = $f
            break
        }
    }
    if ($target -eq $null) {
        Write-Host "Field with instr [$instr] not found"
        return
    }

    $pos = $target.Code.Start - 1
    $target.Delete()

    $r = $doc.Range($pos, $pos)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText>' + $instr + '</w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

Expand-SimpleField $d "m:usercontent zone1"
Expand-SimpleField $d "m:endusercontent"
